$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct a few previously-entered daily figures (column C) ---
$ws.Range("C255").Value = 458
$ws.Range("C258").Value = 553
$ws.Range("C259").Value = 425
$ws.Range("C260").Value = 340
$ws.Range("C261").Value = 262

# Columns L and M are number-formatted as Text ("@"); round-trip the
# number format so the COM layer writes a real number into the cell
# (matching the workbook's existing numeric L/M entries) instead of
# a text string, then restore the original Text format.
$ws.Range("L261").NumberFormat = "general"
$ws.Range("L261").Value = 5
$ws.Range("L261").NumberFormat = "@"

# --- Enter the new day's figures (row 262) ---
$ws.Range("C262").Value = 36
$ws.Range("E262").Value = 35
$ws.Range("F262").Value = 26
$ws.Range("G262").Value = 267

$ws.Range("L262").NumberFormat = "general"
$ws.Range("L262").Value = 1
$ws.Range("L262").NumberFormat = "@"

$ws.Range("M262").NumberFormat = "general"
$ws.Range("M262").Value = 0
$ws.Range("M262").NumberFormat = "@"

# Recalculate so every dependent formula (columns B, H, J, K) picks up
# the new cached values.
$excel.Calculate()

# --- Update the view: scroll the frozen pane down and select the data
#     column that was just filled in ---
$ws.Range("C3:C262").Select()
$excel.ActiveWindow.ScrollRow = 246
$excel.ActiveWindow.ScrollColumn = 2
